$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 9 de Julio de 2020 a las 23:24"

# Row 4: Estados Unidos
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 3206985
$ws.Range("C4").Value = 48053
$ws.Range("D4").Value = 1415302
$ws.Range("E4").Value = 1656111
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 710
$ws.Range("H4").Value = 135572

# Row 5: Brasil
$ws.Range("A5").Value = "Brasil"
$ws.Range("B5").Value = 1755779
$ws.Range("C5").Value = 39583
$ws.Range("D5").Value = 1152467
$ws.Range("E5").Value = 534128
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 1129
$ws.Range("H5").Value = 69184

# Row 19: Alemania
$ws.Range("A19").Value = "Alemania"
$ws.Range("B19").Value = 199162
$ws.Range("C19").Value = 397
$ws.Range("D19").Value = 183600
$ws.Range("E19").Value = 6437
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 10
$ws.Range("H19").Value = 9125

# Row 27: Egipto
$ws.Range("A27").Value = "Egipto"
$ws.Range("B27").Value = 79254
$ws.Range("C27").Value = 950
$ws.Range("D27").Value = 22753
$ws.Range("E27").Value = 52884
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 53
$ws.Range("H27").Value = 3617

# Row 50: Barein
$ws.Range("A50").Value = "Barein"
$ws.Range("B50").Value = 31528
$ws.Range("C50").Value = 597
$ws.Range("D50").Value = 26520
$ws.Range("E50").Value = 4905
$ws.Range("F50").Value = 0
$ws.Range("G50").Value = 5
$ws.Range("H50").Value = 103

# Row 72: Sudan
$ws.Range("A72").Value = "Sudan"
$ws.Range("B72").Value = 10158
$ws.Range("C72").Value = 74
$ws.Range("D72").Value = 5200
$ws.Range("E72").Value = 4317
$ws.Range("F72").Value = 0
$ws.Range("G72").Value = 5
$ws.Range("H72").Value = 641

# Row 105: Somalia
$ws.Range("A105").Value = "Somalia"
$ws.Range("B105").Value = 3038
$ws.Range("C105").Value = 10
$ws.Range("D105").Value = 1209
$ws.Range("E105").Value = 1737
$ws.Range("F105").Value = 0
$ws.Range("G105").Value = 0
$ws.Range("H105").Value = 92

# Row 108: Paraguay
$ws.Range("A108").Value = "Paraguay"
$ws.Range("B108").Value = 2638
$ws.Range("C108").Value = 84
$ws.Range("D108").Value = 1229
$ws.Range("E108").Value = 1389
$ws.Range("F108").Value = 0
$ws.Range("G108").Value = 0
$ws.Range("H108").Value = 20

# Row 128: Yemen
$ws.Range("A128").Value = "Yemen"
$ws.Range("B128").Value = 1356
$ws.Range("C128").Value = 38
$ws.Range("D128").Value = 619
$ws.Range("E128").Value = 376
$ws.Range("F128").Value = 0
$ws.Range("G128").Value = 10
$ws.Range("H128").Value = 361

# Row 132: Ruanda
$ws.Range("A132").Value = "Ruanda"
$ws.Range("B132").Value = 1210
$ws.Range("C132").Value = 16
$ws.Range("D132").Value = 623
$ws.Range("E132").Value = 584
$ws.Range("F132").Value = 0
$ws.Range("G132").Value = 0
$ws.Range("H132").Value = 3

# Row 133: Benin
$ws.Range("A133").Value = "Benin"
$ws.Range("B133").Value = 1199
$ws.Range("C133").Value = 0
$ws.Range("D133").Value = 333
$ws.Range("E133").Value = 845
$ws.Range("F133").Value = 0
$ws.Range("G133").Value = 0
$ws.Range("H133").Value = 21

# Row 151: Togo
$ws.Range("A151").Value = "Togo"
$ws.Range("B151").Value = 704
$ws.Range("C151").Value = 9
$ws.Range("D151").Value = 483
$ws.Range("E151").Value = 206
$ws.Range("F151").Value = 0
$ws.Range("G151").Value = 0
$ws.Range("H151").Value = 15

# Row 152: San Marino
$ws.Range("A152").Value = "San Marino"
$ws.Range("B152").Value = 698
$ws.Range("C152").Value = 0
$ws.Range("D152").Value = 656
$ws.Range("E152").Value = 0
$ws.Range("F152").Value = 0
$ws.Range("G152").Value = 0
$ws.Range("H152").Value = 42

# Row 153: Surinam
$ws.Range("A153").Value = "Surinam"
$ws.Range("B153").Value = 694
$ws.Range("C153").Value = 29
$ws.Range("D153").Value = 456
$ws.Range("E153").Value = 221
$ws.Range("F153").Value = 0
$ws.Range("G153").Value = 0
$ws.Range("H153").Value = 17

# Row 154: Malta
$ws.Range("A154").Value = "Malta"
$ws.Range("B154").Value = 674
$ws.Range("C154").Value = 1
$ws.Range("D154").Value = 656
$ws.Range("E154").Value = 9
$ws.Range("F154").Value = 0
$ws.Range("G154").Value = 0
$ws.Range("H154").Value = 9

# Row 164: Birmania
$ws.Range("A164").Value = "Birmania"
$ws.Range("B164").Value = 319
$ws.Range("C164").Value = 2
$ws.Range("D164").Value = 254
$ws.Range("E164").Value = 59
$ws.Range("F164").Value = 0
$ws.Range("G164").Value = 0
$ws.Range("H164").Value = 6

# Row 184: Seychelles
$ws.Range("A184").Value = "Seychelles"
$ws.Range("B184").Value = 91
$ws.Range("C184").Value = 0
$ws.Range("D184").Value = 11
$ws.Range("E184").Value = 80
$ws.Range("F184").Value = 0
$ws.Range("G184").Value = 0
$ws.Range("H184").Value = 0

# Row 185: Lesoto
$ws.Range("A185").Value = "Lesoto"
$ws.Range("B185").Value = 91
$ws.Range("C185").Value = 0
$ws.Range("D185").Value = 11
$ws.Range("E185").Value = 80
$ws.Range("F185").Value = 0
$ws.Range("G185").Value = 0
$ws.Range("H185").Value = 0

# Row 209: Islas Malvinas
$ws.Range("A209").Value = "Islas Malvinas"
$ws.Range("B209").Value = 13
$ws.Range("C209").Value = 0
$ws.Range("D209").Value = 13
$ws.Range("E209").Value = 0
$ws.Range("F209").Value = 0
$ws.Range("G209").Value = 0
$ws.Range("H209").Value = 0

# Row 210: Groenlandia
$ws.Range("A210").Value = "Groenlandia"
$ws.Range("B210").Value = 13
$ws.Range("C210").Value = 0
$ws.Range("D210").Value = 13
$ws.Range("E210").Value = 0
$ws.Range("F210").Value = 0
$ws.Range("G210").Value = 0
$ws.Range("H210").Value = 0
